# Update Nuuchahnulth data progress workbook.
$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("data")
$graphsSheet = $wb.Worksheets.Item("graphs")

# Current_Page (B4 on the "data" sheet) advances from 583 to 594.
$dataSheet.Range("B4").Value = 594

# Actual pages completed logged for the most recent date (row 25) advances
# from 1176 to 1187 - this feeds the "Actual Percent Remaining" chart series.
$dataSheet.Range("G25").Value = 1187

# Update the active selection / active sheet to match the author's saved view:
# the data sheet's selection moves to G26, and the graphs sheet becomes the
# active (selected) tab.
$dataSheet.Range("G26").Select()
$graphsSheet.Activate()
